$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix header cells with "/" to adopt JSONPointer-style keys
$ws.Range("A1").Value = "/Index"
$ws.Range("B1").Value = "/Array_1()"
$ws.Range("C1").Value = "/Array_2(Int)"

# Move the active selection from E3 to C5
$ws.Range("C5").Select()
